$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.05441666666666667
$ws.Range("H2").Value = 0.16325
$ws.Range("I2").Value = 0.00608027172874025
$ws.Range("J2").Value = 0.006080271728740251
$ws.Range("M2").Value = 0.1124096666666667
$ws.Range("N2").Value = 0.337229
$ws.Range("O2").Value = 0.0063156040276642
$ws.Range("P2").Value = 0.0063960530566531
$ws.Range("Q2").Value = 0.006116959361111111
$ws.Range("R2").Value = 0.05505263425
$ws.Range("S2").Value = 0.00003840058861932469
$ws.Range("T2").Value = 0.00003888974057589052
$ws.Range("G3").Value = 0.05441666666666667
$ws.Range("H3").Value = 0.16325
$ws.Range("I3").Value = 0.00608027172874025
$ws.Range("J3").Value = 0.006080271728740251
$ws.Range("O3").Value = 0.1214504387717248
$ws.Range("P3").Value = 0.1229974910927163
$ws.Range("Q3").Value = 0.1176304586388889
$ws.Range("R3").Value = 1.05867412775
$ws.Range("S3").Value = 0.000738451669306817
$ws.Range("T3").Value = 0.0007478581677970237
$ws.Range("G4").Value = 0.05441666666666667
$ws.Range("H4").Value = 0.16325
$ws.Range("I4").Value = 0.00608027172874025
$ws.Range("J4").Value = 0.006080271728740251
$ws.Range("M4").Value = 10.589294
$ws.Range("N4").Value = 31.767882
$ws.Range("O4").Value = 0.5949469455757395
$ws.Range("P4").Value = 0.6025254612429388
$ws.Range("Q4").Value = 0.5762340818333334
$ws.Range("R4").Value = 5.1861067365
$ws.Range("S4").Value = 0.003617439093284533
$ws.Range("T4").Value = 0.003663518527841621
$ws.Range("G5").Value = 0.05441666666666667
$ws.Range("H5").Value = 0.16325
$ws.Range("I5").Value = 0.00608027172874025
$ws.Range("J5").Value = 0.006080271728740251
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.6716124999999999
$ws.Range("N5").Value = 1.343225
$ws.Range("O5").Value = 0.03773375311758142
$ws.Range("P5").Value = 0.02547627388813791
$ws.Range("Q5").Value = 0.03654691354166666
$ws.Range("R5").Value = 0.21928148125
$ws.Range("S5").Value = 0.0002294314723000946
$ws.Range("T5").Value = 0.0001549026678756884
$ws.Range("G6").Value = 0.05441666666666667
$ws.Range("H6").Value = 0.16325
$ws.Range("I6").Value = 0.00608027172874025
$ws.Range("J6").Value = 0.006080271728740251
$ws.Range("M6").Value = 4.263741333333333
$ws.Range("N6").Value = 12.791224
$ws.Range("O6").Value = 0.2395532585072902
$ws.Range("P6").Value = 0.2426047207195541
$ws.Range("Q6").Value = 0.2320185908888889
$ws.Range("R6").Value = 2.088167318
$ws.Range("S6").Value = 0.001456548905229481
$ws.Range("T6").Value = 0.001475102624650029
$ws.Range("G7").Value = 4.046611333333334
$ws.Range("H7").Value = 12.139834
$ws.Range("I7").Value = 0.4521500120171497
$ws.Range("J7").Value = 0.4521500120171497
$ws.Range("M7").Value = 0.1124096666666667
$ws.Range("N7").Value = 0.337229
$ws.Range("O7").Value = 0.0063156040276642
$ws.Range("P7").Value = 0.0063960530566531
$ws.Range("Q7").Value = 0.4548782311095556
$ws.Range("R7").Value = 4.093904079986
$ws.Range("S7").Value = 0.002855600437003927
$ws.Range("T7").Value = 0.002891975466428026
$ws.Range("G8").Value = 4.046611333333334
$ws.Range("H8").Value = 12.139834
$ws.Range("I8").Value = 0.4521500120171497
$ws.Range("J8").Value = 0.4521500120171497
$ws.Range("O8").Value = 0.1214504387717248
$ws.Range("P8").Value = 0.1229974910927163
$ws.Range("Q8").Value = 8.747407296906447
$ws.Range("R8").Value = 78.72666567215801
$ws.Range("S8").Value = 0.05491381735012346
$ws.Range("T8").Value = 0.05561331707565093
$ws.Range("G9").Value = 4.046611333333334
$ws.Range("H9").Value = 12.139834
$ws.Range("I9").Value = 0.4521500120171497
$ws.Range("J9").Value = 0.4521500120171497
$ws.Range("M9").Value = 10.589294
$ws.Range("N9").Value = 31.767882
$ws.Range("O9").Value = 0.5949469455757395
$ws.Range("P9").Value = 0.6025254612429388
$ws.Range("Q9").Value = 42.85075711239868
$ws.Range("R9").Value = 385.656814011588
$ws.Range("S9").Value = 0.2690052685916371
$ws.Range("T9").Value = 0.2724318945416334
$ws.Range("G10").Value = 4.046611333333334
$ws.Range("H10").Value = 12.139834
$ws.Range("I10").Value = 0.4521500120171497
$ws.Range("J10").Value = 0.4521500120171497
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.6716124999999999
$ws.Range("N10").Value = 1.343225
$ws.Range("O10").Value = 0.03773375311758142
$ws.Range("P10").Value = 0.02547627388813791
$ws.Range("Q10").Value = 2.717754754108333
$ws.Range("R10").Value = 16.30652852465
$ws.Range("S10").Value = 0.0170613169255666
$ws.Range("T10").Value = 0.01151909754467375
$ws.Range("G11").Value = 4.046611333333334
$ws.Range("H11").Value = 12.139834
$ws.Range("I11").Value = 0.4521500120171497
$ws.Range("J11").Value = 0.4521500120171497
$ws.Range("M11").Value = 4.263741333333333
$ws.Range("N11").Value = 12.791224
$ws.Range("O11").Value = 0.2395532585072902
$ws.Range("P11").Value = 0.2426047207195541
$ws.Range("Q11").Value = 17.25370400186845
$ws.Range("R11").Value = 155.283336016816
$ws.Range("S11").Value = 0.1083140087128186
$ws.Range("T11").Value = 0.1096937273887636
$ws.Range("G12").Value = 4.848681666666667
$ws.Range("H12").Value = 14.546045
$ws.Range("I12").Value = 0.5417697162541101
$ws.Range("J12").Value = 0.5417697162541102
$ws.Range("M12").Value = 0.1124096666666667
$ws.Range("N12").Value = 0.337229
$ws.Range("O12").Value = 0.0063156040276642
$ws.Range("P12").Value = 0.0063960530566531
$ws.Range("Q12").Value = 0.5450386899227778
$ws.Range("R12").Value = 4.905348209305
$ws.Range("S12").Value = 0.003421603002040948
$ws.Range("T12").Value = 0.003465187849649184
$ws.Range("G13").Value = 4.848681666666667
$ws.Range("H13").Value = 14.546045
$ws.Range("I13").Value = 0.5417697162541101
$ws.Range("J13").Value = 0.5417697162541102
$ws.Range("O13").Value = 0.1214504387717248
$ws.Range("P13").Value = 0.1229974910927163
$ws.Range("Q13").Value = 10.48121252515722
$ws.Range("R13").Value = 94.33091272641499
$ws.Range("S13").Value = 0.06579816975229451
$ws.Range("T13").Value = 0.06663631584926834
$ws.Range("G14").Value = 4.848681666666667
$ws.Range("H14").Value = 14.546045
$ws.Range("I14").Value = 0.5417697162541101
$ws.Range("J14").Value = 0.5417697162541102
$ws.Range("M14").Value = 10.589294
$ws.Range("N14").Value = 31.767882
$ws.Range("O14").Value = 0.5949469455757395
$ws.Range("P14").Value = 0.6025254612429388
$ws.Range("Q14").Value = 51.34411568074334
$ws.Range("R14").Value = 462.09704112669
$ws.Range("S14").Value = 0.3223242378908178
$ws.Range("T14").Value = 0.3264300481734638
$ws.Range("G15").Value = 4.848681666666667
$ws.Range("H15").Value = 14.546045
$ws.Range("I15").Value = 0.5417697162541101
$ws.Range("J15").Value = 0.5417697162541102
$ws.Range("K15").Value = 2
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.6716124999999999
$ws.Range("N15").Value = 1.343225
$ws.Range("O15").Value = 0.03773375311758142
$ws.Range("P15").Value = 0.02547627388813791
$ws.Range("Q15").Value = 3.256435215854166
$ws.Range("R15").Value = 19.538611295125
$ws.Range("S15").Value = 0.02044300471971473
$ws.Range("T15").Value = 0.01380227367558847
$ws.Range("G16").Value = 4.848681666666667
$ws.Range("H16").Value = 14.546045
$ws.Range("I16").Value = 0.5417697162541101
$ws.Range("J16").Value = 0.5417697162541102
$ws.Range("M16").Value = 4.263741333333333
$ws.Range("N16").Value = 12.791224
$ws.Range("O16").Value = 0.2395532585072902
$ws.Range("P16").Value = 0.2426047207195541
$ws.Range("Q16").Value = 20.67352443434222
$ws.Range("R16").Value = 186.06171990908
$ws.Range("S16").Value = 0.1297827008892421
$ws.Range("T16").Value = 0.1314358907061405
